# Remove the first 10 data rows (the "Distraction\..." video entries),
# which shifts the remaining "Drowsy\..." rows up so they become rows 2-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:11").Delete()
